$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 340, shifting existing rows 340-418 down to 343-421
$ws.Rows("340:342").Insert()

# Row 340
$ws.Cells.Item(340, 1).Value = 6
$ws.Cells.Item(340, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(340, 3).Value = "Metropolitana"
$ws.Cells.Item(340, 4).Value = 44754
$ws.Cells.Item(340, 5).Value = 13
$ws.Cells.Item(340, 6).Value = "Fruta"
$ws.Cells.Item(340, 7).Value = 100107
$ws.Cells.Item(340, 8).Value = "Otros"
$ws.Cells.Item(340, 9).Value = 100107011
$ws.Cells.Item(340, 10).Value = "Tuna"
$ws.Cells.Item(340, 11).Value = "Sin especificar"
$ws.Cells.Item(340, 12).Value = "Especial"
$ws.Cells.Item(340, 13).Value = 70
$ws.Cells.Item(340, 14).Value = 27000
$ws.Cells.Item(340, 15).Value = 27000
$ws.Cells.Item(340, 16).Value = 27000
$ws.Cells.Item(340, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(340, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(340, 19).Value = 1500
$ws.Cells.Item(340, 20).Value = 18

# Row 341
$ws.Cells.Item(341, 1).Value = 6
$ws.Cells.Item(341, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(341, 3).Value = "Metropolitana"
$ws.Cells.Item(341, 4).Value = 44754
$ws.Cells.Item(341, 5).Value = 13
$ws.Cells.Item(341, 6).Value = "Fruta"
$ws.Cells.Item(341, 7).Value = 100107
$ws.Cells.Item(341, 8).Value = "Otros"
$ws.Cells.Item(341, 9).Value = 100107011
$ws.Cells.Item(341, 10).Value = "Tuna"
$ws.Cells.Item(341, 11).Value = "Sin especificar"
$ws.Cells.Item(341, 12).Value = "Primera"
$ws.Cells.Item(341, 13).Value = 80
$ws.Cells.Item(341, 14).Value = 25000
$ws.Cells.Item(341, 15).Value = 25000
$ws.Cells.Item(341, 16).Value = 25000
$ws.Cells.Item(341, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(341, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(341, 19).Value = 1389
$ws.Cells.Item(341, 20).Value = 18

# Row 342
$ws.Cells.Item(342, 1).Value = 6
$ws.Cells.Item(342, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(342, 3).Value = "Metropolitana"
$ws.Cells.Item(342, 4).Value = 44754
$ws.Cells.Item(342, 5).Value = 13
$ws.Cells.Item(342, 6).Value = "Fruta"
$ws.Cells.Item(342, 7).Value = 100107
$ws.Cells.Item(342, 8).Value = "Otros"
$ws.Cells.Item(342, 9).Value = 100107011
$ws.Cells.Item(342, 10).Value = "Tuna"
$ws.Cells.Item(342, 11).Value = "Sin especificar"
$ws.Cells.Item(342, 12).Value = "Segunda"
$ws.Cells.Item(342, 13).Value = 50
$ws.Cells.Item(342, 14).Value = 22000
$ws.Cells.Item(342, 15).Value = 22000
$ws.Cells.Item(342, 16).Value = 22000
$ws.Cells.Item(342, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(342, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(342, 19).Value = 1222
$ws.Cells.Item(342, 20).Value = 18
